$d = $word.ActiveDocument

# Locate the title paragraph ("Dheeraj Chand") so we can insert the missing
# contact-info line directly underneath it.
$searchRange = $d.Content
$searchRange.Find.Execute("Dheeraj Chand", $true, $false, $false, $false, $false,
                           $true, 1, $false, "", 0) | Out-Null

$titlePara = $searchRange.Paragraphs(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)   # wdCollapseEnd - sit right after the title paragraph mark

# Re-insert the title paragraph (unchanged) together with a brand-new, cleanly
# formatted paragraph for the contact info, via raw OOXML. Using InsertXML lets
# us add the paragraph without it inheriting the bold/28pt run formatting of the
# title line.
$contactText = "202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX"

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr><w:jc w:val="center"/></w:pPr>
            <w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t>Dheeraj Chand</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr><w:jc w:val="center"/></w:pPr>
            <w:r><w:t>$contactText</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$titleRange.InsertXML($xml)
